$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append two new AVC test-run rows (13 & 14) to the results table, copying
# the formatting of the last existing row (12) so number/date formats and
# column styles stay identical.
# ---------------------------------------------------------------------------

# Row 13 - commit 4d722633c6f4821e59a9ca2b393870674f751c87
$ws.Range("A12:F12").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)
$ws.Range("A13").Value = 219
$ws.Range("B13").Value = 0.00000126721215565
$ws.Range("C13").Value = 15671770293
$ws.Range("F13").Value = "4d722633c6f4821e59a9ca2b393870674f751c87"
$ws.Range("D13").Value = "Last AVC Test for Today"
$ws.Range("E13").Value = 42869

# Row 14 - commit 63b114717e2642c74648886f53259cba73b21231
$ws.Range("A12:F12").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)
$ws.Range("A14").Value = 239
$ws.Range("B14").Value = 0.00000126721215565
$ws.Range("C14").Value = 15998233425
$ws.Range("F14").Value = "63b114717e2642c74648886f53259cba73b21231"
$ws.Range("D14").Value = "some more c (useless tmp variables) deleted"
$ws.Range("E14").Value = 42870

$excel.CutCopyMode = $false

# Selection moves to D17 (matches the author's saved cursor position).
$ws.Range("D17").Select()

# ---------------------------------------------------------------------------
# Extend both charts' series so they plot through row 14 instead of row 12.
# ---------------------------------------------------------------------------
$charts = $ws.ChartObjects()

$co1 = $charts.Item(1)
$chart1 = $co1.Chart
$series1 = $chart1.SeriesCollection().Item(1)
$series1.Formula = "=SERIES(Tabelle1!`$C`$1,,Tabelle1!`$C`$2:`$C`$14,1)"

$co2 = $charts.Item(2)
$chart2 = $co2.Chart
$series2 = $chart2.SeriesCollection().Item(1)
$series2.Formula = "=SERIES(Tabelle1!`$B`$1,,Tabelle1!`$B`$2:`$B`$14,1)"

# ---------------------------------------------------------------------------
# The two new rows pushed the chart anchors further down the sheet; move
# both chart objects to their new position (size is unchanged).
# ---------------------------------------------------------------------------
$co1.Left = 3.5293700787401576
$co1.Top = 332.294094488189
$co1.Width = 365.831654158465
$co1.Height = 216

$co2.Left = 402.89055179625984
$co2.Top = 333.1763779527559
$co2.Width = 367.60595617984245
$co2.Height = 216
